# logboek.xlsx update:
# - Add a new logboek entry in row 19 (date 29/12/2024, 1.5 "Gewerkte uren")
#   describing the registration-system fix / role-based URL routing work.
# - Move the active selection to A20 (just below the newly filled row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = "registratie: als je niet alles invuld krijg je een error als het email adres al in database staat krijg je een error en het verplaten van urls met role ipv alles ingelogt alles door getest"
$ws.Range("B19").Value = 45655
$ws.Range("D19").Value = 1.5

$ws.Range("A20").Select() | Out-Null
